# Slide 1: adjust bar-chart gridline/bar/label positions and update
# several value labels + the figure caption N count.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)  # top-level group shape holding the whole figure

$shp = $grp.GroupItems.Item("pl5")
$shp.Top = 349.7398681640625

$shp = $grp.GroupItems.Item("pl6")
$shp.Top = 298.5958557128906

$shp = $grp.GroupItems.Item("pl7")
$shp.Top = 247.45181

$shp = $grp.GroupItems.Item("pl8")
$shp.Top = 196.3078

$shp = $grp.GroupItems.Item("pl9")
$shp.Top = 145.1638

$shp = $grp.GroupItems.Item("rc10")
$shp.Top = 379.5739
$shp.Height = 21.310001373291016

$shp = $grp.GroupItems.Item("rc11")
$shp.Top = 397.3322
$shp.Height = 3.5517

$shp = $grp.GroupItems.Item("rc12")
$shp.Top = 276.5755310058594
$shp.Height = 124.30835

$shp = $grp.GroupItems.Item("rc14")
$shp.Top = 308.5405
$shp.Height = 92.34331

$shp = $grp.GroupItems.Item("rc15")
$shp.Top = 365.3672
$shp.Height = 35.51661682128906

$shp = $grp.GroupItems.Item("rc16")
$shp.Top = 397.3322
$shp.Height = 3.5517

$shp = $grp.GroupItems.Item("tx17")
$shp.Top = 333.3543
$shp.Height = 10.39882
$shp.TextFrame.TextRange.Text = "6"

$shp = $grp.GroupItems.Item("tx18")
$shp.Top = 350.888
$shp.TextFrame.TextRange.Text = "(4%)"

$shp = $grp.GroupItems.Item("tx19")
$shp.Top = 351.28622

$shp = $grp.GroupItems.Item("tx20")
$shp.Top = 368.6463

$shp = $grp.GroupItems.Item("tx21")
$shp.Top = 230.3489
$shp.TextFrame.TextRange.Text = "35"

$shp = $grp.GroupItems.Item("tx22")
$shp.Top = 247.88961791992188
$shp.TextFrame.TextRange.Text = "(24%)"

$shp = $grp.GroupItems.Item("tx23")
$shp.TextFrame.TextRange.Text = "65"

$shp = $grp.GroupItems.Item("tx24")
$shp.TextFrame.TextRange.Text = "(45%)"

$shp = $grp.GroupItems.Item("tx25")
$shp.Top = 262.3209
$shp.Height = 10.39882
$shp.TextFrame.TextRange.Text = "26"

$shp = $grp.GroupItems.Item("tx26")
$shp.Top = 279.85467529296875

$shp = $grp.GroupItems.Item("tx27")
$shp.Top = 319.1476

$shp = $grp.GroupItems.Item("tx28")
$shp.Top = 336.6813

$shp = $grp.GroupItems.Item("tx29")
$shp.Top = 351.28622

$shp = $grp.GroupItems.Item("tx30")
$shp.Top = 368.6463

$shp = $grp.GroupItems.Item("tx33")
$shp.Top = 345.2633361816406

$shp = $grp.GroupItems.Item("tx34")
$shp.Top = 294.1193

$shp = $grp.GroupItems.Item("tx35")
$shp.Top = 242.9694

$shp = $grp.GroupItems.Item("tx36")
$shp.Top = 191.8313

$shp = $grp.GroupItems.Item("tx37")
$shp.Top = 140.687255859375

$shp = $grp.GroupItems.Item("pl39")
$shp.Top = 349.7398681640625

$shp = $grp.GroupItems.Item("pl40")
$shp.Top = 298.5958557128906

$shp = $grp.GroupItems.Item("pl41")
$shp.Top = 247.45181

$shp = $grp.GroupItems.Item("pl42")
$shp.Top = 196.3078

$shp = $grp.GroupItems.Item("pl43")
$shp.Top = 145.1638

$shp = $grp.GroupItems.Item("tx61")
$shp.TextFrame.TextRange.Text = "Level of involvement in treatment choice (N=144)."
